$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$nRows = 24
$nCols = 13
$data = New-Object 'object[,]' $nRows,$nCols

$data[0,0] = 7.885082333931324
$data[0,1] = 7.007767235305031
$data[0,2] = 11.98118886863211
$data[0,3] = 34.71544382563388
$data[0,4] = 3.678963699735785
$data[0,5] = 0
$data[0,6] = 27.29803617798842
$data[0,7] = 9.524122541916787
$data[0,8] = 0
$data[0,9] = 9.583718933443111
$data[0,10] = 60.46959249696625
$data[0,11] = 0
$data[0,12] = 27.53840627938159
$data[1,0] = 7.968954733943305
$data[1,1] = 7.048721818118671
$data[1,2] = 11.98947848744773
$data[1,3] = 35.18446773958331
$data[1,4] = 3.682851571810195
$data[1,5] = 0
$data[1,6] = 27.6381568321496
$data[1,7] = 9.55804152759163
$data[1,8] = 0
$data[1,9] = 9.593210467867189
$data[1,10] = 57.13609579565174
$data[1,11] = 0
$data[1,12] = 27.8551623561513
$data[2,0] = 8.022725822857376
$data[2,1] = 7.075166537152021
$data[2,2] = 11.99893571001962
$data[2,3] = 35.48775941154377
$data[2,4] = 3.685340911345668
$data[2,5] = 0
$data[2,6] = 27.85784480651443
$data[2,7] = 9.581694126228106
$data[2,8] = 0
$data[2,9] = 9.601779054518268
$data[2,10] = 54.97768428405774
$data[2,11] = 0
$data[2,12] = 28.06099155482351
$data[3,0] = 8.045210293847555
$data[3,1] = 7.086269990229058
$data[3,2] = 12.00387590002504
$data[3,3] = 35.61517421618255
$data[3,4] = 3.686381175719935
$data[3,5] = 0
$data[3,6] = 27.95008264254452
$data[3,7] = 9.592036207694063
$data[3,8] = 0
$data[3,9] = 9.605954415504565
$data[3,10] = 54.07039149097371
$data[3,11] = 0
$data[3,12] = 28.14769556508875
$data[4,0] = 8.04897841288728
$data[4,1] = 7.088133476700607
$data[4,2] = 12.00476147242302
$data[4,3] = 35.63656132505193
$data[4,4] = 3.686555475630685
$data[4,5] = 0
$data[4,6] = 27.96556209498961
$data[4,7] = 9.593795782563545
$data[4,8] = 0
$data[4,9] = 9.606688854745856
$data[4,10] = 53.91807385305455
$data[4,11] = 0
$data[4,12] = 28.16226274744211
$data[5,0] = 8.023026737949614
$data[5,1] = 7.075314957416251
$data[5,2] = 11.99899795263378
$data[5,3] = 35.48946233742306
$data[5,4] = 3.685354835903653
$data[5,5] = 0
$data[5,6] = 27.85907779252945
$data[5,7] = 9.581830764249045
$data[5,8] = 0
$data[5,9] = 9.601832604429914
$data[5,10] = 54.96555993642838
$data[5,11] = 0
$data[5,12] = 28.06214946063062
$data[6,0] = 7.913530171277305
$data[6,1] = 7.021619186707994
$data[6,2] = 11.9831335661963
$data[6,3] = 34.87396811441744
$data[6,4] = 3.680283122322213
$data[6,5] = 0
$data[6,6] = 27.41304757633114
$data[6,7] = 9.535226820382176
$data[6,8] = 0
$data[6,9] = 9.586419405023936
$data[6,10] = 59.34355298329994
$data[6,11] = 0
$data[6,12] = 27.64525511255687
$data[7,0] = 7.716813183728945
$data[7,1] = 6.926601802015522
$data[7,2] = 11.98720440298697
$data[7,3] = 33.7898243639255
$data[7,4] = 3.671141257296307
$data[7,5] = 0
$data[7,6] = 26.62524251442751
$data[7,7] = 9.46659683443106
$data[7,8] = 0
$data[7,9] = 9.578194193051543
$data[7,10] = 67.03203838171132
$data[7,11] = 0
$data[7,12] = 26.91885700645192
$data[8,0] = 7.583224564207361
$data[8,1] = 6.863029011277801
$data[8,2] = 12.01235605307332
$data[8,3] = 33.07022106769214
$data[8,4] = 3.664904677901908
$data[8,5] = 0
$data[8,6] = 26.10050100678478
$data[8,7] = 9.430528346877622
$data[8,8] = 0
$data[8,9] = 9.58591478027879
$data[8,10] = 72.1242302232743
$data[8,11] = 0
$data[8,12] = 26.442412328337
$data[9,0] = 7.524822473836869
$data[9,1] = 6.835456693024399
$data[9,2] = 12.02875802910499
$data[9,3] = 32.76005546295936
$data[9,4] = 3.662169464860057
$data[9,5] = 0
$data[9,6] = 25.87379313402954
$data[9,7] = 9.417344678408572
$data[9,8] = 0
$data[9,9] = 9.592491171494673
$data[9,10] = 74.31896086743527
$data[9,11] = 0
$data[9,12] = 26.23849560247476
$data[10,0] = 7.503047769016325
$data[10,1] = 6.825209209800573
$data[10,2] = 12.03569384477428
$data[10,3] = 32.64512221070988
$data[10,4] = 3.661148173449325
$data[10,5] = 0
$data[10,6] = 25.789696274398
$data[10,7] = 9.412825194701806
$data[10,8] = 0
$data[10,9] = 9.595428166069544
$data[10,10] = 75.13250889096126
$data[10,11] = 0
$data[10,12] = 26.16315797849509
$data[11,0] = 7.5077221710913
$data[11,1] = 6.827407582458908
$data[11,2] = 12.0341676807419
$data[11,3] = 32.66976226165723
$data[11,4] = 3.661367485829329
$data[11,5] = 0
$data[11,6] = 25.80772964165696
$data[11,7] = 9.413777364366076
$data[11,8] = 0
$data[11,9] = 9.594775670878544
$data[11,10] = 74.95807761253921
$data[11,11] = 0
$data[11,12] = 26.17929898995575
$data[12,0] = 7.523024220645097
$data[12,1] = 6.834609748115803
$data[12,2] = 12.02931405353312
$data[12,3] = 32.75054904388578
$data[12,4] = 3.662085153502175
$data[12,5] = 0
$data[12,6] = 25.86683914981736
$data[12,7] = 9.416963326485631
$data[12,8] = 0
$data[12,9] = 9.592723813758406
$data[12,10] = 74.38624379358112
$data[12,11] = 0
$data[12,12] = 26.23225960397723
$data[13,0] = 7.532441578286281
$data[13,1] = 6.839046484811928
$data[13,2] = 12.02643578391527
$data[13,3] = 32.80036293208867
$data[13,4] = 3.66252662554848
$data[13,5] = 0
$data[13,6] = 25.90327449512836
$data[13,7] = 9.418976683299437
$data[13,8] = 0
$data[13,9] = 9.591525338314767
$data[13,10] = 74.03369294728516
$data[13,11] = 0
$data[13,12] = 26.26494565897204
$data[14,0] = 7.587088956062707
$data[14,1] = 6.864858020802044
$data[14,2] = 12.01138494653799
$data[14,3] = 33.09084098788187
$data[14,4] = 3.665085465343083
$data[14,5] = 0
$data[14,6] = 26.11556061629106
$data[14,7] = 9.431455593537061
$data[14,8] = 0
$data[14,9] = 9.585547173602096
$data[14,10] = 71.97834533294952
$data[14,11] = 0
$data[14,12] = 26.4559999077526
$data[15,0] = 7.621220058734298
$data[15,1] = 6.881037482157762
$data[15,2] = 12.00343046013724
$data[15,3] = 33.27347472013484
$data[15,4] = 3.666681194534433
$data[15,5] = 0
$data[15,6] = 26.24888315283892
$data[15,7] = 9.439943275067071
$data[15,8] = 0
$data[15,9] = 9.582668884660251
$data[15,10] = 70.68623549382708
$data[15,11] = 0
$data[15,12] = 26.57651512734666
$data[16,0] = 7.641074221214107
$data[16,1] = 6.890470298471612
$data[16,2] = 11.99932147835349
$data[16,3] = 33.38013526061768
$data[16,4] = 3.667608611696912
$data[16,5] = 0
$data[16,6] = 26.32669447125075
$data[16,7] = 9.445128059348487
$data[16,8] = 0
$data[16,9] = 9.581301433039044
$data[16,10] = 69.93159027746167
$data[16,11] = 0
$data[16,12] = 26.64703660726524
$data[17,0] = 7.647834767217494
$data[17,1] = 6.893685875432168
$data[17,2] = 11.99800996849127
$data[17,3] = 33.41652455779959
$data[17,4] = 3.667924272417202
$data[17,5] = 0
$data[17,6] = 26.35323294862177
$data[17,7] = 9.44693528215309
$data[17,8] = 0
$data[17,9] = 9.58088772859335
$data[17,10] = 69.67411436898649
$data[17,11] = 0
$data[17,12] = 26.6711197023326
$data[18,0] = 7.61756367799285
$data[18,1] = 6.879302027188014
$data[18,2] = 12.00422889201177
$data[18,3] = 33.25386562014395
$data[18,4] = 3.666510334510635
$data[18,5] = 0
$data[18,6] = 26.23457383830399
$data[18,7] = 9.439008332163899
$data[18,8] = 0
$data[18,9] = 9.582945426227727
$data[18,10] = 70.82496905342188
$data[18,11] = 0
$data[18,12] = 26.56356116971089
$data[19,0] = 7.518520381605441
$data[19,1] = 6.832489045557653
$data[19,2] = 12.0307199224692
$data[19,3] = 32.72675121613076
$data[19,4] = 3.661873965654603
$data[19,5] = 0
$data[19,6] = 25.8494294660902
$data[19,7] = 9.416014623033753
$data[19,8] = 0
$data[19,9] = 9.593314323855545
$data[19,10] = 74.55468170697964
$data[19,11] = 0
$data[19,12] = 26.21665241964713
$data[20,0] = 7.455777238555806
$data[20,1] = 6.803022204803779
$data[20,2] = 12.05226278398478
$data[20,3] = 32.39696388670961
$data[20,4] = 3.658928116330366
$data[20,5] = 0
$data[20,6] = 25.60794439349139
$data[20,7] = 9.403747725384504
$data[20,8] = 0
$data[20,9] = 9.6026972397355
$data[20,10] = 76.88997784431471
$data[20,11] = 0
$data[20,12] = 26.0009119386947
$data[21,0] = 7.489082400250535
$data[21,1] = 6.818646030019433
$data[21,2] = 12.04037425994054
$data[21,3] = 32.57161430784607
$data[21,4] = 3.660492716736544
$data[21,5] = 0
$data[21,6] = 25.7358840546525
$data[21,7] = 9.41003905552339
$data[21,8] = 0
$data[21,9] = 9.597448946158698
$data[21,10] = 75.65295193947983
$data[21,11] = 0
$data[21,12] = 26.11503851783797
$data[22,0] = 7.619216005943506
$data[22,1] = 6.880086218241571
$data[22,2] = 12.00386647608298
$data[22,3] = 33.26272572118816
$data[22,4] = 3.666587549132208
$data[22,5] = 0
$data[22,6] = 26.24103946075127
$data[22,7] = 9.439430070221398
$data[22,8] = 0
$data[22,9] = 9.582819507319426
$data[22,10] = 70.76228434332216
$data[22,11] = 0
$data[22,12] = 26.56941380792842
$data[23,0] = 7.768107203251182
$data[23,1] = 6.951209304507169
$data[23,2] = 11.98226209070112
$data[23,3] = 34.06976437213633
$data[23,4] = 3.67352934931339
$data[23,5] = 0
$data[23,6] = 26.82895750561045
$data[23,7] = 9.482677924713256
$data[23,8] = 0
$data[23,9] = 9.578027676522181
$data[23,10] = 65.04934424225102
$data[23,11] = 0
$data[23,12] = 27.10544383035086

$range = $ws.Range("C2:O25")
$range.Value = $data

Write-Output "Updated loading_percent values for rows 2-25"